# Updated cryptos list values (Price and Volume(1h) columns) on the
# "Sheet1" worksheet, cells D2:D51 and E2:E51 (rows 8 and 33 unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "586.32", "1.00") but must
# stay plain text cells, matching the original inline-string cells. Setting
# Range.Value directly on such strings makes Excel auto-convert them to
# numbers. Routing the text through a text-producing formula on a scratch
# cell, then Copy + PasteSpecial(xlPasteValues) into the target, preserves
# the text type without touching any NumberFormat/Style (which would leave
# an unwanted style-table diff).
$stage = $ws.Range("ZZ1")
$stage.Formula = "=""65.476.71"""
$stage.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$stage.Formula = "=""3.424.56"""
$stage.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$stage.Formula = "=""586.32"""
$stage.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$stage.Formula = "=""137.87"""
$stage.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$stage.Formula = "=""3.426.30"""
$stage.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$stage.Formula = "=""0.497"""
$stage.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$stage.Formula = "=""7.23"""
$stage.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$stage.Formula = "=""4.003.38"""
$stage.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$stage.Formula = "=""26.26"""
$stage.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$stage.Formula = "=""3.428.62"""
$stage.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$stage.Formula = "=""65.394.15"""
$stage.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$stage.Formula = "=""9.79"""
$stage.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$stage.Formula = "=""5.91"""
$stage.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$stage.Formula = "=""13.63"""
$stage.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$stage.Formula = "=""389.76"""
$stage.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$stage.Formula = "=""0.557"""
$stage.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$stage.Formula = "=""73.03"""
$stage.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$stage.Formula = "=""3.561.21"""
$stage.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$stage.Formula = "=""1.00"""
$stage.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$stage.Formula = "=""7.13"""
$stage.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$stage.Formula = "=""8.17"""
$stage.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$stage.Formula = "=""2.21"""
$stage.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$stage.Formula = "=""3.427.16"""
$stage.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$stage.Formula = "=""0.144"""
$stage.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$stage.Formula = "=""22.92"""
$stage.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$stage.Formula = "=""172.25"""
$stage.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$stage.Formula = "=""6.82"""
$stage.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$stage.Formula = "=""1.14"""
$stage.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$stage.Formula = "=""1.46"""
$stage.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$stage.Formula = "=""4.74"""
$stage.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$stage.Formula = "=""0.0764"""
$stage.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$stage.Formula = "=""0.818"""
$stage.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$stage.Formula = "=""43.51"""
$stage.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$stage.Formula = "=""0.999"""
$stage.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$stage.Formula = "=""4.41"""
$stage.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$stage.Formula = "=""1.61"""
$stage.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$stage.Formula = "=""1.10"""
$stage.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$stage.Formula = "=""22.39"""
$stage.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$stage.Formula = "=""6.54"""
$stage.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$stage.Formula = "=""2.07"""
$stage.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$stage.Formula = "=""2.183.98"""
$stage.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$stage.Clear()

# Column E ("Volume(1h)") values are padded percentage strings (e.g.
# "  -0.77%  ") which Excel always keeps as text, so a direct assignment is fine.
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("E7").Value = "  -2.82%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  -6.82%  "
$ws.Range("E11").Value = "  -9.55%  "
$ws.Range("E12").Value = "  -7.25%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  -10.30%  "
$ws.Range("E15").Value = "  -8.52%  "
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -10.60%  "
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("E22").Value = "  -6.15%  "
$ws.Range("E23").Value = "  -6.76%  "
$ws.Range("E24").Value = "  -5.75%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("E27").Value = "  -8.89%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -8.69%  "
$ws.Range("E30").Value = "  -9.13%  "
$ws.Range("E31").Value = "  -9.41%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("E34").Value = "  -7.12%  "
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  -9.63%  "
$ws.Range("E38").Value = "  -10.91%  "
$ws.Range("E39").Value = "  -7.50%  "
$ws.Range("E40").Value = "  -9.96%  "
$ws.Range("E41").Value = "  -6.74%  "
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -13.05%  "
$ws.Range("E46").Value = "  -10.57%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("E49").Value = "  -7.68%  "
$ws.Range("E50").Value = "  -14.67%  "
$ws.Range("E51").Value = "  -7.09%  "
